# Adipoq-Adipor2.xlsx -- "update scripts wuth new tpm"
#
# The NATMI LR-pair sheet was re-run with a newer TPM matrix. The
# sending/target cluster labels are unchanged (FAPs, Adipoq, Adipor2,
# ECs, MuSCs, Resolving-Mac) but:
#   - all the numeric statistics for the 4 existing "FAPs -> Adipoq/
#     Adipor2" rows were recomputed against the new TPM values, and
#   - 4 brand-new rows were added for the "MuSCs -> Adipoq/Adipor2"
#     sending cluster (rows 6-9), extending the used range from
#     A1:T5 to A1:T9.
#
# Columns (1-indexed) for every data row:
#   A Sending cluster        K Receptor-expressing cells
#   B Ligand symbol          L Receptor detection rate
#   C Receptor symbol        M Receptor average expression value
#   D Target cluster         N Receptor total expression value
#   E Ligand-expressing cells   O Receptor derived specificity (avg)
#   F Ligand detection rate     P Receptor derived specificity (total)
#   G Ligand average expression value  Q Edge average expression weight
#   H Ligand total expression value    R Edge total expression weight
#   I Ligand derived specificity (avg) S Edge average expression derived specificity
#   J Ligand derived specificity (tot) T Edge total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is: row number, then the 20 values for columns A..T
# in order. Text columns (A-D) are plain strings; the rest are numbers.
$rows = @(
    @(2, "FAPs",  "Adipoq", "Adipor2", "ECs",           1, 0.3333333333333333, 0.1884053333333333, 0.5652160000000001, 0.603254837001998,  0.603254837001998,  3, 1, 15.232775,         45.698325,  0.4461934273123386, 0.4461934273123386, 2.869936051466667,  25.8294244632,       0.2691683432646677,  0.2691683432646677),
    @(3, "FAPs",  "Adipoq", "Adipor2", "FAPs",          1, 0.3333333333333333, 0.1884053333333333, 0.5652160000000001, 0.603254837001998,  0.603254837001998,  3, 1, 5.529931333333333,  16.589794,  0.1619809269435077, 0.1619809269435077, 1.041868556167111,  9.376817005504,      0.09771577768073827, 0.09771577768073827),
    @(4, "FAPs",  "Adipoq", "Adipor2", "MuSCs",         1, 0.3333333333333333, 0.1884053333333333, 0.5652160000000001, 0.603254837001998,  0.603254837001998,  3, 1, 5.511105333333333,  16.533316,  0.1614294819531772, 0.1614294819531772, 1.038321637361778,  9.344894736256,      0.09738311582298088, 0.09738311582298088),
    @(5, "FAPs",  "Adipoq", "Adipor2", "Resolving-Mac", 1, 0.3333333333333333, 0.1884053333333333, 0.5652160000000001, 0.603254837001998,  0.603254837001998,  3, 1, 7.865586333333333,  23.596759,  0.2303961637909764, 0.2303961637909764, 1.481918414993778,  13.337265734944,     0.1389876002336111,  0.1389876002336111),
    @(6, "MuSCs", "Adipoq", "Adipor2", "ECs",           1, 0.3333333333333333, 0.1239093333333333, 0.371728,            0.396745162998002,  0.396745162998002,  3, 1, 15.232775,         45.698325,  0.4461934273123386, 0.4461934273123386, 1.887482995066667,  16.9873469556,       0.1770250840476709,  0.1770250840476709),
    @(7, "MuSCs", "Adipoq", "Adipor2", "FAPs",          1, 0.3333333333333333, 0.1239093333333333, 0.371728,            0.396745162998002,  0.396745162998002,  3, 1, 5.529931333333333,  16.589794,  0.1619809269435077, 0.1619809269435077, 0.6852101048924444, 6.166890944031999,   0.06426514926276941, 0.06426514926276941),
    @(8, "MuSCs", "Adipoq", "Adipor2", "MuSCs",         1, 0.3333333333333333, 0.1239093333333333, 0.371728,            0.396745162998002,  0.396745162998002,  3, 1, 5.511105333333333,  16.533316,  0.1614294819531772, 0.1614294819531772, 0.682877387783111,  6.145896490048,      0.06404636613019631, 0.06404636613019631),
    @(9, "MuSCs", "Adipoq", "Adipor2", "Resolving-Mac", 1, 0.3333333333333333, 0.1239093333333333, 0.371728,            0.396745162998002,  0.396745162998002,  3, 1, 7.865586333333333,  23.596759,  0.2303961637909764, 0.2303961637909764, 0.974619558839111,  8.771576029552,      0.0914085635573653,  0.0914085635573653)
)

foreach ($row in $rows) {
    $r = $row[0]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }
}
